# Add the new users to the roster sheet (Лист1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("Саша", "slunka322", "Ксюша🫧", "Joe", "Andris", "Родион", "f¡dgy")
$ids   = @(931591593, 886229823, 1260079637, 5450675821, 635469686, 1030349543, 1122159904)
$nicks = @("@yaderon", "@slunka322", "@kssyusshh", "@SpecCorvo", "@Tut_dedus", "@Nx1dxr", "@fiidgy")

# Fill column A (names) first, then column C (nicknames), then column B (ids) -
# matches the order the strings were interned into the shared-strings table.
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $names[$i]
}

for ($i = 0; $i -lt $nicks.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $nicks[$i]
}

for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $ids[$i]
}

# Stamp the new data rows (name + id columns) with their own cell style so
# they pick up a dedicated cellXfs entry, same as the source workbook.
$ws.Range("A2:B8").Locked = $true

$ws.Range("C8").Select()
